$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.575.03'
$ws.Range('E2').Value = '  -5.38%  '

# Row 3
$ws.Range('D3').Value = '2.983.28'
$ws.Range('E3').Value = '  -7.16%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '543.58'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -5.56%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '152.77'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -8.73%  '

# Row 7
$ws.Range('E7').Value = '  -0.06%  '

# Row 8
$ws.Range('E8').Value = '  -6.18%  '

# Row 9
$ws.Range('D9').Value = '2.984.06'
$ws.Range('E9').Value = '  -6.94%  '

# Row 10
$ws.Range('E10').Value = '  -7.00%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.19'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -8.23%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.364'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -7.21%  '

# Row 13
$ws.Range('D13').Value = '3.496.55'
$ws.Range('E13').Value = '  -7.32%  '

# Row 14
$ws.Range('E14').Value = '  -3.76%  '

# Row 15
$ws.Range('D15').Value = '61.668.48'
$ws.Range('E15').Value = '  -5.32%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '23.55'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -8.16%  '

# Row 17
$ws.Range('D17').Value = '2.985.22'
$ws.Range('E17').Value = '  -6.66%  '

# Row 18
$ws.Range('E18').Value = '  -7.40%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '388.02'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -5.92%  '

# Row 20
$ws.Range('E20').Value = '  -4.96%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.84'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -8.32%  '

# Row 22
$ws.Range('E22').Value = '  -8.14%  '

# Row 23
$ws.Range('E23').Value = '  -0.12%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '64.89'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -6.98%  '

# Row 25
$ws.Range('E25').Value = '  -5.27%  '

# Row 26
$ws.Range('E26').Value = '  -8.48%  '

# Row 27
$ws.Range('E27').Value = '  -0.52%  '

# Row 28
$ws.Range('D28').Value = '0.0₃0936'
$ws.Range('E28').Value = '  -11.18%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.30'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -7.02%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = "Normal"

# Row 31
$ws.Range('E31').Value = '  -7.60%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '20.24'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -6.59%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '158.61'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.81%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.98'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -6.82%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.60'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -7.90%  '

# Row 36
$ws.Range('E36').Value = '  -6.83%  '

# Row 37
$ws.Range('E37').Value = '  -7.05%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.58'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -9.16%  '

# Row 39
$ws.Range('D39').Value = '2.441.67'
$ws.Range('E39').Value = '  -11.08%  '

# Row 40
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '37.20'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -5.02%  '

# Row 41
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.86'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -7.20%  '

# Row 42
$ws.Range('E42').Value = '  -8.66%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.658'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -7.96%  '

# Row 44
$ws.Range('E44').Value = '  -6.71%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.998'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.23%  '

# Row 46
$ws.Range('E46').Value = '  -6.97%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.90'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -13.08%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0955'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.82%  '

# Row 49
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '10.49'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.14%  '

# Row 50
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '19.63'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -8.84%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '263.71'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -11.21%  '
